# Refresh the cryptos list: update each coin's Price (column D) and
# Volume(1h) change (column E) cell to the latest scraped reading.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.568.99"
$ws.Range("E2").Value = "  -2.65%  "
$ws.Range("D3").Value = "1.756.46"
$ws.Range("E3").Value = "  -3.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4496"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3632"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.105"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.15%  "
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.056"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.232"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.01%  "
$ws.Range("D16").Value = "1.757.90"
$ws.Range("E16").Value = "  -3.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06431"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.869"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.06%  "
$ws.Range("D23").Value = "27.618.91"
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("E24").Value = "  -2.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.096"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("D28").Value = "1.959.57"
$ws.Range("E28").Value = "  -3.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.135"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.095"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09092"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.555"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.75%  "
$ws.Range("E34").Value = "  +4.63%  "
$ws.Range("E35").Value = "  -6.38%  "
$ws.Range("E36").Value = "  -1.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6422"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2102"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.943"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.194"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.398"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.848"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5918"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.715"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.972"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.163"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06876"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.71%  "
